$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$m = $s.Master
$cs = $m.ColorScheme
$cs.Colors(5).RGB = 1111111
